$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.805.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.14%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.624.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.08%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.15'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.58%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.59%  '

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.396'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.62%  '

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.82'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.96%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.78%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.89'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.101.44'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.736.86'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.03%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +13.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.572.05'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.43%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.35%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.78'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.93%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '347.53'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.74%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.38%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.42'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.68%  '

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.83%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.10%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.92%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.28'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.75%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '548.76'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.29%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.60%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.92%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.36'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.59%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.11'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.35%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.11'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.17%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.417'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.93%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.97'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.98'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '168.18'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.31%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.86%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.38'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.77%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0585'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.27%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.72%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.637'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.86%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0252'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.48%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0969'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.29'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.13%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0230'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +17.64%  '
